{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst oldServerText = \"Ensure that the transfer on the server side has terminated with an appropriate error message. \";\nconst newServerText = \"Check the server log to ensure that no RRQ was received by the server\";\nconst code6Text = \"Check the server log to ensure that an Error packet with code 6 was received by the Server\";\n\n// Find the target \"server side\" paragraph that still carries the trailing\n// space (there are similarly-worded paragraphs earlier in the doc that do\n// NOT have the trailing space / the following \"code 6\" sibling \u2014 skip those).\nlet targetIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text === oldServerText) {\n    const next = paragraphs.items[i + 2];\n    if (next && next.text === code6Text) {\n      targetIndex = i;\n      break;\n    }\n  }\n}\n\nif (targetIndex === -1) {\n  throw new Error(\"Could not locate target paragraph for edit\");\n}\n\n// 1) Rewrite the \"server side\" bullet text.\nparagraphs.items[targetIndex].insertText(newServerText, \"Replace\");\n\n// 2) Delete the \"Check the server log ... code 6\" bullet (two after target,\n// i.e. right after the unchanged \"client side\" bullet) and the blank\n// paragraph that immediately follows it.\nconst code6Paragraph = paragraphs.items[targetIndex + 2];\nconst blankParagraph = paragraphs.items[targetIndex + 3];\ncode6Paragraph.load(\"text\");\nblankParagraph.load(\"text\");\nawait context.sync();\n\nif (code6Paragraph.text !== code6Text) {\n  throw new Error(\"Unexpected paragraph while deleting code-6 bullet\");\n}\nif (blankParagraph.text !== \"\") {\n  throw new Error(\"Unexpected paragraph while deleting blank paragraph\");\n}\n\nblankParagraph.delete();\ncode6Paragraph.delete();\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Get-ParaText($para) {\n    return $para.Range.Text.TrimEnd([char]13)\n}\n\n$oldServerText = \"Ensure that the transfer on the server side has terminated with an appropriate error message. \"\n$newServerText = \"Check the server log to ensure that no RRQ was received by the server\"\n$code6Text = \"Check the server log to ensure that an Error packet with code 6 was received by the Server\"\n\n# Locate the \"server side\" bullet that still has the trailing space AND is\n# followed two paragraphs later by the \"code 6\" bullet (there are similarly\n# worded paragraphs earlier in the document that must not be touched).\n$n = $d.Paragraphs.Count\n$targetIndex = -1\nfor ($i = 1; $i -le $n; $i++) {\n    if ((Get-ParaText $d.Paragraphs($i)) -eq $oldServerText) {\n        if (($i + 2) -le $n -and (Get-ParaText $d.Paragraphs($i + 2)) -eq $code6Text) {\n            $targetIndex = $i\n            break\n        }\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not locate target paragraph for edit\"\n}\n\n$code6Index = $targetIndex + 2\n$blankIndex = $targetIndex + 3\n\nif ((Get-ParaText $d.Paragraphs($blankIndex)) -ne \"\") {\n    throw \"Unexpected paragraph while deleting blank paragraph\"\n}\n\n# 1) Rewrite the \"server side\" bullet text (keeps the paragraph mark/formatting).\n$d.Paragraphs($targetIndex).Range.Text = $newServerText\n\n# 2) Delete the \"Check the server log ... code 6\" bullet and the blank\n# paragraph right after it. Delete from the highest index down so earlier\n# indices stay valid.\n$d.Paragraphs($blankIndex).Range.Delete()\n$d.Paragraphs($code6Index).Range.Delete()\n"}
